$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (date) values between rows 2/3 and rows 4/5.
# Row 2 and 3 move from 2022-10-04 (44838) to 2022-10-12 (44846)
# Row 4 and 5 move from 2022-10-12 (44846) to 2022-10-04 (44838)
$ws.Range("D2").Value = 44846
$ws.Range("D3").Value = 44846
$ws.Range("D4").Value = 44838
$ws.Range("D5").Value = 44838
